# update game page ui
#
# The "walk" town action (a leisurely stroll) was repurposed into a
# "game" town action (an arcade / game hall). Rename the two localization
# keys and update the accompanying simplified-Chinese title/description,
# then move the sheet's scroll position & selection to match where the
# author was working (a couple of columns to the left of where it was).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 holds the localization keys used by the game code.
$ws.Range("AE1").Value = "game_title"
$ws.Range("AF1").Value = "game_desc"

# Row 2 holds the simplified-Chinese (the base/default language) strings.
$ws.Range("AE2").Value = "游戏厅"
$ws.Range("AF2").Value = "放松你的心情\再来亿把"

# The workbook was left scrolled to column W with AF3 selected (previously
# it was scrolled to column Z with AJ3 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 23
$ws.Range("AF3").Select()
